# Append run: 2026-01-09 06:30 JST
# The scraper re-ran and produced a fresh (shorter) result set for the
# "ランサーズ" sheet. Net effect vs. the previous snapshot:
#   - only 4 listings remain (rows 2-5) instead of 14 (rows 2-15)
#   - timestamps bumped to the new run time
#   - three listings are carried over from the previous run (now ranked
#     differently / re-scored) and one listing is brand new
#   - a few column widths were tweaked

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Drop the hyperlinks collection up front; row deletes below would
# otherwise leave stale hyperlink entries pointing at ranges that no
# longer hold URLs. We rebuild the ones we need at the end. ---
$ws.Hyperlinks.Delete()

# --- Trim the sheet down from 14 data rows to the 4 that remain. ---
# Keep (by original row number): 3 -> new 2, 4 -> new 3, 5 -> new 4 (rewritten),
# 11 -> new 5. Delete everything else.
$ws.Rows("12:15").Delete()
$ws.Rows("6:10").Delete()
$ws.Rows("2:2").Delete()

# --- Refresh the "fetched at" timestamp on every surviving row. ---
$ws.Range("A2:A5").Value = "2026-01-09 06:30:45"

# --- Row 4 becomes an entirely new listing (same price bucket as the old
# "Ecommerce" row it replaced, but a different job, URL, score, and no
# skill tags this time). ---
$ws.Range("B4").Value = "【急募】大手保険システム会社でのPJ推進支援(PM・PL経験者募集/都内常駐)"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5467981"
$ws.Range("G4").Value = 40
$ws.Range("H4").ClearContents()

# --- Column width tweaks. COM ColumnWidth uses the "character" unit which
# Excel stores in the XML `width` attribute offset by the default 5/6
# padding factor for this workbook's font -- subtract that back out so the
# saved width lands on the exact target values (B:41, D:30, H:12). ---
$pad = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 41 - $pad
$ws.Columns.Item(4).ColumnWidth = 30 - $pad
$ws.Columns.Item(8).ColumnWidth = 12 - $pad

# --- Rebuild the hyperlinks for the URL column on the surviving rows. ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5467745")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5467910")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5467981")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5467882")
